$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new date header in column U (treat as text, matching how the other
# date headers in row 1 are stored, not an auto-converted Excel date)
$ws.Range("U1").NumberFormat = "@"
$ws.Range("U1").Value = "2025-06-23"
$ws.Range("U1").Font.Bold = $true
$ws.Range("U1").Borders.LineStyle = 1
$ws.Range("U1").HorizontalAlignment = -4108
$ws.Range("U1").VerticalAlignment = -4160

# Update running totals in column S for the two student rows
$ws.Range("S2").Value = 16
$ws.Range("S3").Value = 16

# Mark the new date's attendance as absent for both students
$ws.Range("U2").Value = "❌"
$ws.Range("U3").Value = "❌"
